# "fixed bug with ip" - summary_report.xlsx
#
# The IP-address cell held a stringified Python list (e.g. "['0.0.0.0', '0.0.0.0']")
# instead of a single address, and the pass/error counters next to it were stale.
# This updates the values and narrows column B back down now that it no longer
# needs to hold the long list-like string.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B2: IP address - was "['0.0.0.0', '0.0.0.0']", now a plain address
$ws.Range("B2").Value = "0.0.0.0"

# The counters in B3:B6 are stored as text (shared strings), not numbers, in the
# original workbook. Assigning a numeric-looking string directly would make Excel
# store it as a real number (or mark it with a "number stored as text" quote
# prefix), so format the cell as text first, write the value, then drop the
# number-format override we just added so the cell keeps the plain default style.
function Set-TextValue($range, [string]$text) {
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.ClearFormats()
}

Set-TextValue $ws.Range("B3") "1632"   # Conf passed: 1634 -> 1632
Set-TextValue $ws.Range("B4") "4"      # Conf error: 2 -> 4
Set-TextValue $ws.Range("B5") "69"     # HW passed: 71 -> 69
Set-TextValue $ws.Range("B6") "2"      # HW error: 0 -> 2

# Column B was sized to fit the old long list-like string; shrink it back down
# (Excel's ColumnWidth is pixel-quantized, so this lands on the closest
# representable width to the authored 7.7109375 character width).
$ws.Columns.Item(2).ColumnWidth = 6.86
